$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 983.2941
$ws.Range("J17").Value = 1066.6207
$ws.Range("L17").Value = 3199.8621
$ws.Range("N17").Value = -3535.8621

# Row 129
$ws.Range("H129").Value = 1856.138
$ws.Range("J129").Value = 2297.8096
$ws.Range("L129").Value = 6893.4288
$ws.Range("N129").Value = -16893.4288

# Row 132
$ws.Range("H132").Value = 4791.145
$ws.Range("I132").Value = 3634.362
$ws.Range("J132").Value = 10890.546
$ws.Range("K132").Value = 10903.086
$ws.Range("L132").Value = 32671.638
$ws.Range("M132").Value = -8373.085999999999
$ws.Range("N132").Value = -37731.638

# Row 137
$ws.Range("H137").Value = 2117.2683
$ws.Range("I137").Value = 2955.8667
$ws.Range("J137").Value = 1633.4615
$ws.Range("K137").Value = 8867.6001
$ws.Range("L137").Value = 4900.3845
$ws.Range("M137").Value = -6317.6001
$ws.Range("N137").Value = -10000.3845

# Row 138
$ws.Range("H138").Value = 2756.3647
$ws.Range("I138").Value = 1747.3214
$ws.Range("J138").Value = 3370.5652
$ws.Range("K138").Value = 5241.9642
$ws.Range("L138").Value = 10111.6956
$ws.Range("M138").Value = -101.9642000000003
$ws.Range("N138").Value = -20391.6956


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 16
$ws.Range("H16").Value = 2812.2856
$ws.Range("I16").Value = 1937.2
$ws.Range("K16").Value = 1937.2
$ws.Range("M16").Value = -1650.2

# Row 43
$ws.Range("H43").Value = 18059
$ws.Range("J43").Value = 18059
$ws.Range("L43").Value = 18059
$ws.Range("N43").Value = -18685

# Row 110
$ws.Range("H110").Value = 1316.6364
$ws.Range("I110").Value = 899.3077
$ws.Range("J110").Value = 1919.4445
$ws.Range("K110").Value = 899.3077
$ws.Range("L110").Value = 1919.4445
$ws.Range("M110").Value = 1145.6923
$ws.Range("N110").Value = -6009.4445

# Row 132
$ws.Range("H132").Value = 4668.7383
$ws.Range("I132").Value = 3901.6316
$ws.Range("J132").Value = 5748.3706
$ws.Range("K132").Value = 11704.8948
$ws.Range("L132").Value = 17245.1118
$ws.Range("M132").Value = -9174.8948
$ws.Range("N132").Value = -22305.1118


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 10
$ws.Range("H10").Value = 18167.834
$ws.Range("I10").Value = 3666.6667
$ws.Range("J10").Value = 32669
$ws.Range("K10").Value = 3666.6667
$ws.Range("L10").Value = 32669
$ws.Range("M10").Value = -3526.6667
$ws.Range("N10").Value = -32949

# Row 17
$ws.Range("H17").Value = 750
$ws.Range("J17").Value = 750
$ws.Range("L17").Value = 750
$ws.Range("N17").Value = -1094

# Row 24
$ws.Range("H24").Value = 1945
$ws.Range("I24").Value = 1945
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 1945
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -1710
$ws.Range("N24").ClearContents()

# Row 86
$ws.Range("H86").Value = 1076.5
$ws.Range("I86").Value = 944.5
$ws.Range("J86").Value = 1340.5
$ws.Range("K86").Value = 944.5
$ws.Range("L86").Value = 1340.5
$ws.Range("M86").Value = 178.5
$ws.Range("N86").Value = -3586.5

# Row 89
$ws.Range("H89").Value = 1076.5
$ws.Range("I89").Value = 944.5
$ws.Range("J89").Value = 1340.5
$ws.Range("K89").Value = 4722.5
$ws.Range("L89").Value = 6702.5
$ws.Range("M89").Value = 893.5
$ws.Range("N89").Value = -17934.5


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 1173.25
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# Row 9
$ws.Range("H9").Value = 71201.664
$ws.Range("J9").Value = 71201.664
$ws.Range("L9").Value = 71201.664
$ws.Range("N9").Value = -71537.664

# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# Row 19
$ws.Range("H19").Value = 957.3333
$ws.Range("I19").Value = 538.9
$ws.Range("J19").Value = 3049.5
$ws.Range("K19").Value = 538.9
$ws.Range("L19").Value = 3049.5
$ws.Range("M19").Value = -368.9
$ws.Range("N19").Value = -3389.5

# Row 24
$ws.Range("H24").Value = 957.3333
$ws.Range("I24").Value = 538.9
$ws.Range("J24").Value = 3049.5
$ws.Range("K24").Value = 538.9
$ws.Range("L24").Value = 3049.5
$ws.Range("M24").Value = -368.9
$ws.Range("N24").Value = -3389.5

# Row 31
$ws.Range("H31").Value = 2197.4546
$ws.Range("I31").Value = 1136.8793
$ws.Range("J31").Value = 3697.7805
$ws.Range("K31").Value = 1136.8793
$ws.Range("L31").Value = 3697.7805
$ws.Range("M31").Value = -841.8793000000001
$ws.Range("N31").Value = -4287.7805

# Row 33
$ws.Range("H33").Value = 4653
$ws.Range("I33").Value = 4653
$ws.Range("K33").Value = 4653
$ws.Range("M33").Value = -4274

# Row 34
$ws.Range("H34").Value = 2197.4546
$ws.Range("I34").Value = 1136.8793
$ws.Range("J34").Value = 3697.7805
$ws.Range("K34").Value = 1136.8793
$ws.Range("L34").Value = 3697.7805
$ws.Range("M34").Value = -934.8793000000001
$ws.Range("N34").Value = -4101.7805

# Row 36
$ws.Range("H36").Value = 2679
$ws.Range("I36").Value = 1518.5
$ws.Range("J36").Value = 5000
$ws.Range("K36").Value = 1518.5
$ws.Range("L36").Value = 5000
$ws.Range("M36").Value = -1130.5
$ws.Range("N36").Value = -5776

# Row 40
$ws.Range("H40").Value = 2679
$ws.Range("I40").Value = 1518.5
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 1518.5
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -1358.5
$ws.Range("N40").Value = -5320


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 1616.375
$ws.Range("I9").Value = 820.5
$ws.Range("K9").Value = 820.5
$ws.Range("M9").Value = -650.5

# Row 13
$ws.Range("H13").Value = 451
$ws.Range("I13").Value = 268.33334
$ws.Range("J13").Value = 999
$ws.Range("K13").Value = 268.33334
$ws.Range("L13").Value = 999
$ws.Range("M13").Value = -129.33334
$ws.Range("N13").Value = -1277

# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

# Row 23
$ws.Range("H23").Value = 414
$ws.Range("J23").Value = 414
$ws.Range("L23").Value = 414
$ws.Range("N23").Value = -860

# Row 43
$ws.Range("H43").Value = 4510
$ws.Range("I43").Value = 4510
$ws.Range("K43").Value = 4510
$ws.Range("M43").Value = -4359

# Row 46
$ws.Range("H46").Value = 7695.4
$ws.Range("I46").Value = 7695.4
$ws.Range("K46").Value = 7695.4
$ws.Range("M46").Value = -7539.4

# Row 70
$ws.Range("H70").Value = 10050
$ws.Range("I70").Value = 12325
$ws.Range("J70").Value = 5500
$ws.Range("K70").Value = 12325
$ws.Range("L70").Value = 5500
$ws.Range("M70").Value = -12055
$ws.Range("N70").Value = -6040

# Row 73
$ws.Range("H73").Value = 10050
$ws.Range("I73").Value = 12325
$ws.Range("J73").Value = 5500
$ws.Range("K73").Value = 12325
$ws.Range("L73").Value = 5500
$ws.Range("M73").Value = -11389
$ws.Range("N73").Value = -7372

# Row 80
$ws.Range("H80").Value = 2923.923
$ws.Range("I80").Value = 2775.625
$ws.Range("J80").Value = 3161.2
$ws.Range("K80").Value = 2775.625
$ws.Range("L80").Value = 3161.2
$ws.Range("M80").Value = -1777.625
$ws.Range("N80").Value = -5157.2

# Row 83
$ws.Range("H83").Value = 2923.923
$ws.Range("I83").Value = 2775.625
$ws.Range("J83").Value = 3161.2
$ws.Range("K83").Value = 13878.125
$ws.Range("L83").Value = 15806
$ws.Range("M83").Value = -8886.125
$ws.Range("N83").Value = -25790

# Row 122
$ws.Range("H122").Value = 3235.2666
$ws.Range("I122").Value = 3202.3333
$ws.Range("J122").Value = 3284.6667
$ws.Range("K122").Value = 9606.999899999999
$ws.Range("L122").Value = 9854.000100000001
$ws.Range("M122").Value = -7156.999899999999
$ws.Range("N122").Value = -14754.0001

# Row 123
$ws.Range("H123").Value = 20149.818
$ws.Range("J123").Value = 20149.818
$ws.Range("L123").Value = 20149.818
$ws.Range("N123").Value = -25049.818


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4448.3335
$ws.Range("I7").Value = 4500
$ws.Range("J7").Value = 4438
$ws.Range("K7").Value = 4500
$ws.Range("L7").Value = 4438
$ws.Range("M7").Value = -4388
$ws.Range("N7").Value = -4662

# Row 9
$ws.Range("H9").Value = 382.42856
$ws.Range("I9").Value = 296.16666
$ws.Range("J9").Value = 900
$ws.Range("K9").Value = 296.16666
$ws.Range("L9").Value = 900
$ws.Range("M9").Value = -72.16665999999998
$ws.Range("N9").Value = -1348

# Row 40
$ws.Range("H40").Value = 43482180
$ws.Range("I40").Value = 76926630
$ws.Range("J40").Value = 4392.9
$ws.Range("K40").Value = 76926630
$ws.Range("L40").Value = 4392.9
$ws.Range("M40").Value = -76926494
$ws.Range("N40").Value = -4664.9

# Row 82
$ws.Range("H82").Value = 1645.1578
$ws.Range("I82").Value = 1450.8889
$ws.Range("J82").Value = 1820
$ws.Range("K82").Value = 1450.8889
$ws.Range("L82").Value = 1820
$ws.Range("M82").Value = -1089.8889
$ws.Range("N82").Value = -2542

# Row 85
$ws.Range("H85").Value = 1645.1578
$ws.Range("I85").Value = 1450.8889
$ws.Range("J85").Value = 1820
$ws.Range("K85").Value = 1450.8889
$ws.Range("L85").Value = 1820
$ws.Range("M85").Value = -202.8888999999999
$ws.Range("N85").Value = -4316

# Row 126
$ws.Range("H126").Value = 4448.3335
$ws.Range("I126").Value = 4500
$ws.Range("J126").Value = 4438
$ws.Range("K126").Value = 13500
$ws.Range("L126").Value = 13314
$ws.Range("M126").Value = -11030
$ws.Range("N126").Value = -18254


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 19267.5
$ws.Range("I51").Value = 13535
$ws.Range("J51").Value = 25000
$ws.Range("K51").Value = 13535
$ws.Range("L51").Value = 25000
$ws.Range("M51").Value = -13025
$ws.Range("N51").Value = -26020

# Row 126
$ws.Range("H126").Value = 3063.5
$ws.Range("I126").Value = 2702.1538
$ws.Range("J126").Value = 4629.3335
$ws.Range("K126").Value = 8106.4614
$ws.Range("L126").Value = 13888.0005
$ws.Range("M126").Value = -5636.4614
$ws.Range("N126").Value = -18828.0005

